$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Log Campeonato Torneo 03 2022")

$ws.Cells.Item(7, 1).Value = 'Heredia, Santiago'
$ws.Cells.Item(7, 3).Value = 'semifinal'
$ws.Cells.Item(7, 4).Value = 500.0
$ws.Cells.Item(8, 1).Value = 'Perot, Alejandro'
$ws.Cells.Item(8, 3).Value = 'semifinal'
$ws.Cells.Item(8, 4).Value = 500.0
$ws.Cells.Item(9, 1).Value = 'Lofeudo, Juan Jose'
$ws.Cells.Item(9, 3).Value = 'cuartos'
$ws.Cells.Item(9, 4).Value = 250.0
$ws.Cells.Item(10, 1).Value = 'Mendez, Carlos'
$ws.Cells.Item(10, 3).Value = 'cuartos'
$ws.Cells.Item(10, 4).Value = 250.0
$ws.Cells.Item(11, 1).Value = 'Nohara, Andres'
$ws.Cells.Item(11, 3).Value = 'cuartos'
$ws.Cells.Item(11, 4).Value = 250.0
$ws.Cells.Item(12, 1).Value = 'Velazquez, Marcelo'
$ws.Cells.Item(12, 3).Value = 'cuartos'
$ws.Cells.Item(12, 4).Value = 250.0
$ws.Cells.Item(13, 1).Value = 'Palamedi, Cristian'
$ws.Cells.Item(13, 3).Value = 'octavos'
$ws.Cells.Item(13, 4).Value = 150.0
$ws.Cells.Item(15, 1).Value = 'Palamedi, Cristian'
$ws.Cells.Item(15, 3).Value = 'primero'
$ws.Cells.Item(15, 4).Value = 500.0
$ws.Cells.Item(16, 1).Value = 'Bonelli, Marcos'
$ws.Cells.Item(16, 3).Value = 'segundo'
$ws.Cells.Item(16, 4).Value = 375.0
$ws.Cells.Item(17, 1).Value = 'Escudero, Martin'
$ws.Cells.Item(17, 3).Value = 'semifinal'
$ws.Cells.Item(17, 4).Value = 250.0
$ws.Cells.Item(18, 1).Value = 'Reniero, Gian'
$ws.Cells.Item(18, 3).Value = 'semifinal'
$ws.Cells.Item(18, 4).Value = 250.0
$ws.Cells.Item(20, 1).Value = 'Levin, Raul'
$ws.Cells.Item(20, 3).Value = 'cuartos'
$ws.Cells.Item(20, 4).Value = 125.0
$ws.Cells.Item(21, 1).Value = 'Vetrano, Luca'
$ws.Cells.Item(21, 3).Value = 'cuartos'
$ws.Cells.Item(21, 4).Value = 125.0
$ws.Cells.Item(22, 1).Value = 'Vigouroux, Jorge'
$ws.Cells.Item(22, 3).Value = 'cuartos'
$ws.Cells.Item(22, 4).Value = 125.0
$ws.Cells.Item(23, 1).Value = 'Boetti, Genaro'
$ws.Cells.Item(23, 3).Value = 'octavos'
$ws.Cells.Item(23, 4).Value = 75.0
$ws.Cells.Item(24, 1).Value = 'Colavini, Daniel'
$ws.Cells.Item(24, 3).Value = 'octavos'
$ws.Cells.Item(24, 4).Value = 75.0
$ws.Cells.Item(25, 1).Value = 'Garcia, Francisco'
$ws.Cells.Item(25, 3).Value = 'octavos'
$ws.Cells.Item(25, 4).Value = 75.0
$ws.Cells.Item(26, 1).Value = 'Vigouroux, Manuel'
$ws.Cells.Item(26, 3).Value = 'octavos'
$ws.Cells.Item(26, 4).Value = 75.0
$ws.Cells.Item(28, 1).Value = 'Prettis, Juan'
$ws.Cells.Item(28, 3).Value = 'primero'
$ws.Cells.Item(28, 4).Value = 250.0
$ws.Cells.Item(29, 1).Value = 'Bonelli, Marcos'
$ws.Cells.Item(29, 3).Value = 'segundo'
$ws.Cells.Item(29, 4).Value = 190.0
$ws.Cells.Item(30, 1).Value = 'Alvarez, Sebastian'
$ws.Cells.Item(30, 3).Value = 'semifinal'
$ws.Cells.Item(30, 4).Value = 125.0
$ws.Cells.Item(31, 1).Value = 'Suarez, Milton'
$ws.Cells.Item(31, 3).Value = 'semifinal'
$ws.Cells.Item(31, 4).Value = 125.0
$ws.Cells.Item(32, 1).Value = 'Colavini, Daniel'
$ws.Cells.Item(32, 3).Value = 'cuartos'
$ws.Cells.Item(32, 4).Value = 65.0
$ws.Cells.Item(33, 1).Value = 'Escudero, Martin'
$ws.Cells.Item(33, 3).Value = 'cuartos'
$ws.Cells.Item(33, 4).Value = 65.0
$ws.Cells.Item(34, 1).Value = 'Reniero, Gian'
$ws.Cells.Item(34, 3).Value = 'cuartos'
$ws.Cells.Item(34, 4).Value = 65.0
$ws.Cells.Item(35, 1).Value = 'Savino, Leandro'
$ws.Cells.Item(35, 3).Value = 'cuartos'
$ws.Cells.Item(35, 4).Value = 65.0
$ws.Cells.Item(36, 1).Value = 'Boetti, Genaro'
$ws.Cells.Item(36, 3).Value = 'octavos'
$ws.Cells.Item(36, 4).Value = 40.0
$ws.Cells.Item(37, 1).Value = 'Depaoli, Luciano'
$ws.Cells.Item(37, 3).Value = 'octavos'
$ws.Cells.Item(37, 4).Value = 40.0
$ws.Cells.Item(38, 1).Value = 'Larrosa, Jorge'
$ws.Cells.Item(38, 3).Value = 'octavos'
$ws.Cells.Item(38, 4).Value = 40.0
$ws.Cells.Item(39, 1).Value = 'Morello, Manuel'
$ws.Cells.Item(39, 3).Value = 'octavos'
$ws.Cells.Item(39, 4).Value = 40.0
$ws.Cells.Item(40, 1).Value = 'Musuruana, Francisco'
$ws.Cells.Item(40, 3).Value = 'octavos'
$ws.Cells.Item(40, 4).Value = 40.0
$ws.Cells.Item(41, 1).Value = 'Pilotti, Paz'
$ws.Cells.Item(41, 3).Value = 'octavos'
$ws.Cells.Item(41, 4).Value = 40.0
$ws.Cells.Item(42, 1).Value = 'Vetrano, Luca'
$ws.Cells.Item(42, 3).Value = 'octavos'
$ws.Cells.Item(42, 4).Value = 40.0
$ws.Cells.Item(43, 1).Value = 'Vigouroux, Manuel'
$ws.Cells.Item(43, 3).Value = 'octavos'
$ws.Cells.Item(43, 4).Value = 40.0
$ws.Cells.Item(44, 1).Value = 'Escalante, Samuel'
$ws.Cells.Item(44, 3).Value = '16avos'
$ws.Cells.Item(44, 4).Value = 30.0
$ws.Cells.Item(45, 1).Value = 'Palamedi, Uma'
$ws.Cells.Item(45, 3).Value = '16avos'
$ws.Cells.Item(45, 4).Value = 30.0
